$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Asturias / Gipuzkoa ordering: Asturias now ranks above Gipuzkoa/Guipuzcoa
# and receives freshly updated figures; Gipuzkoa keeps its previous figures but
# moves down one row.
$ws.Range("A23").Value = "Asturias"
$ws.Range("B23").Value = 1827
$ws.Range("C23").Value = 414
$ws.Range("D23").Value = 1285
$ws.Range("E23").Value = 128

$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B24").Value = 1803
$ws.Range("C24").Value = 4680
$ws.Range("D24").Value = 4694
$ws.Range("E24").Value = 107

# Murcia updated figures
$ws.Range("B31").Value = 1413
$ws.Range("C31").Value = 275
$ws.Range("D31").Value = 1044
$ws.Range("E31").Value = 94

# Melilla updated figures
$ws.Range("B54").Value = 98
$ws.Range("C54").Value = 18
$ws.Range("D54").Value = 78

# Ceuta updated figures
$ws.Range("B55").Value = 93
$ws.Range("C55").Value = 15

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 21:52"
